# Generate Report for Handoff
# Adds a new source file (f985ac45-cc97-4923-968f-27f64848b7e5.md) to the
# localization-status report: one new row on the "Overview" sheet and on
# each language sheet ("zh-cn", "de-de"), inserted just above the existing
# ".localization-config" row so that row keeps trailing the table.

$wb = $excel.ActiveWorkbook

$newFileId   = "f985ac45-cc97-4923-968f-27f64848b7e5"
$newFileMd   = "$newFileId.md"
$srcCommit   = "e939a4f65d4a59ae908292f88ac1e696bb80e10f"
$newMdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/$newFileMd"

$xlfHash     = "ddc61e80f922add0bf30f95beae8572f3d11b5d6"
$zhXlfName   = "$newFileId.$xlfHash.zh-cn.xlf"
$deXlfName   = "$newFileId.$xlfHash.de-de.xlf"
$zhXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1bcfa9706bc44677592a9bab3de3a4ea63218c76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$zhXlfName"
$deXlfUrl    = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1cb4b1accb4db6d2b0758f686859a0edb7100d2a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$deXlfName"

$zhHandoffDt = "2016-03-08 20:39:30"
$deHandoffDt = "2016-03-08 20:39:37"

$configUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/.localization-config"
$oldMdUrl    = "https://github.com/OpenLocalizationTest/oltest/blob/$srcCommit/e2e/d8194462-b772-4b80-88da-27074f322671.md"
$oldZhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0abfa9706bc44677592a9bab3de3a4ea63218c76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d8194462-b772-4b80-88da-27074f322671.ad4faf2954d36368c827993ab9c28e26e41b718a.zh-cn.xlf"
$oldDeXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0cb4b1accb4db6d2b0758f686859a0edb7100d2a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d8194462-b772-4b80-88da-27074f322671.ad4faf2954d36368c827993ab9c28e26e41b718a.de-de.xlf"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value2 = $newFileMd
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = "Ready for handoff"

$ws.Hyperlinks.Add($ws.Range("A2"), $oldMdUrl, "", "", "d8194462-b772-4b80-88da-27074f322671.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $newMdUrl, "", "", $newFileMd)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value2 = $newFileMd
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = $zhXlfName
$ws.Range("D3").Value2 = $zhHandoffDt
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Range("D4").Value2 = "0001-01-01 00:00:00"

$ws.Hyperlinks.Add($ws.Range("A2"), $oldMdUrl, "", "", "d8194462-b772-4b80-88da-27074f322671.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $oldZhXlfUrl, "", "", "d8194462-b772-4b80-88da-27074f322671.ad4faf2954d36368c827993ab9c28e26e41b718a.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $newMdUrl, "", "", $newFileMd)
$ws.Hyperlinks.Add($ws.Range("C3"), $zhXlfUrl, "", "", $zhXlfName)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Hyperlinks.Delete()
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value2 = $newFileMd
$ws.Range("B3").Value2 = "Ready for handoff"
$ws.Range("C3").Value2 = $deXlfName
$ws.Range("D3").Value2 = $deHandoffDt
$ws.Range("G3").Value2 = "0001-01-01 00:00:00"
$ws.Range("H3").Value2 = "Include"

$ws.Range("D4").Value2 = "0001-01-01 00:00:00"

$ws.Hyperlinks.Add($ws.Range("A2"), $oldMdUrl, "", "", "d8194462-b772-4b80-88da-27074f322671.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $oldDeXlfUrl, "", "", "d8194462-b772-4b80-88da-27074f322671.ad4faf2954d36368c827993ab9c28e26e41b718a.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $newMdUrl, "", "", $newFileMd)
$ws.Hyperlinks.Add($ws.Range("C3"), $deXlfUrl, "", "", $deXlfName)
$ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", ".localization-config")
